$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B's width for the new "Próg podobieństwa" (similarity
# threshold) field — column B already exists, just give it an explicit
# width like its siblings (target ~19.85546875 chars).
$ws.Columns("B:B").ColumnWidth = 19

# Header row - swap A3/B3 so L.p. moves to A.
$ws.Range("A3").Value = "L.p."

# Row 4 - new sample record (C4 first so "wszystko_JP" lands in the shared
# string table ahead of "Próg podobieństwa").
$ws.Range("C4").Value = "wszystko_JP"

# New "Próg podobieństwa" header goes into B3.
$ws.Range("B3").Value = "Próg podobieństwa"

$ws.Range("B4").Value = 0.85
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0

# Row 5 - new sample record
$ws.Range("C5").Value = "wszystko_JP2"
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1

# Update selection to match the saved view state.
$ws.Range("E14").Select()
